# Größere Zeilen können verglichen werden
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Remove the "Abkürzungen" header in A1, keep the cell/style empty.
$ws.Range("A1").Value = $null

# B2 becomes the new "ADF=Applieba" text (was "adshdaskjhd").
$ws.Range("B2").Value = "ADF=Applieba"

# C2 ("hallo") and D2 ("hi ") are dropped entirely.
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = $null

# Row 4 ("27.05. Tegel-Munich 1200-1400") is removed completely.
$ws.Rows.Item(4).Delete()

# Column widths change slightly for columns A and B (nearest value the
# pixel-quantised ColumnWidth setter can represent).
$ws.Columns.Item(1).ColumnWidth = 28.3333333333333
$ws.Columns.Item(2).ColumnWidth = 25.3333333333333

# Update the stored selection on Sheet1 to D2.
$ws.Range("D2").Select()

# Sheet3's single column gets a bit wider too (whole declared range, cols 1-1025).
$ws3.Columns.Item(1).ColumnWidth = 14.0
